$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 3.65
$ws.Range("I2").Value = 4.4
$ws.Range("J2").Value = 2.18
$ws.Range("K2").Value = 2.25
$ws.Range("L2").Value = 4.55
$ws.Range("M2").Value = 1.27
$ws.Range("N2").Value = 3.1
$ws.Range("O2").Value = 1.8
$ws.Range("P2").Value = 1.8
$ws.Range("Q2").Value = 2.85
$ws.Range("V2").Value = 1.83
$ws.Range("X2").Value = 7.9
$ws.Range("Y2").Value = 8.25
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 13.5
$ws.Range("AD2").Value = 7.1
$ws.Range("AH2").Value = 12.5
$ws.Range("AI2").Value = 25
$ws.Range("AJ2").Value = 14.5
$ws.Range("AK2").Value = 75
$ws.Range("AL2").Value = 45
$ws.Range("AM2").Value = 45

# Row 4 updates
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 4.33
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.91
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 9.5
$ws.Range("AG4").Value = 251
$ws.Range("AM4").Value = 41
$ws.Range("AO4").Value = 9.5

# Row 5 updates
$ws.Range("G5").Value = 3
$ws.Range("I5").Value = 2.35
$ws.Range("L5").Value = 3.2
$ws.Range("O5").Value = 2.35
$ws.Range("P5").Value = 1.57
$ws.Range("Q5").Value = 4.5
$ws.Range("R5").Value = 1.18
$ws.Range("Y5").Value = 12
$ws.Range("AL5").Value = 21
$ws.Range("AO5").Value = 7.5
$ws.Range("AP5").Value = 1.8
$ws.Range("AQ5").Value = 2.05
